$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 203, shifting existing rows 203:311 down to 204:312.
$ws.Rows("203:203").Insert()

# Populate the newly inserted row 203 with the new record's data.
$ws.Range("A203").Value = 4
$ws.Range("B203").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C203").Value = "Los Lagos"
$ws.Range("D203").Value = 44719
$ws.Range("E203").Value = 10
$ws.Range("F203").Value = 100112045
$ws.Range("G203").Value = "Zapallo"
$ws.Range("H203").Value = "Paine"
$ws.Range("I203").Value = "1a (guarda)"
$ws.Range("J203").Value = 800
$ws.Range("K203").Value = 500
$ws.Range("L203").Value = 500
$ws.Range("M203").Value = 500
$ws.Range("N203").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O203").Value = "Región de O'Higgins"
$ws.Range("P203").Value = 500
$ws.Range("Q203").Value = 1
$ws.Range("R203").Value = "Hortaliza"
